$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1730.3182
$ws.Range("I80").Value = 913.44446
$ws.Range("J80").Value = 2295.8462
$ws.Range("K80").Value = 2740.33338
$ws.Range("L80").Value = 6887.5386
$ws.Range("M80").Value = -1742.33338
$ws.Range("N80").Value = -8883.5386

$ws.Range("H83").Value = 1730.3182
$ws.Range("I83").Value = 913.44446
$ws.Range("J83").Value = 2295.8462
$ws.Range("K83").Value = 8221.00014
$ws.Range("L83").Value = 20662.6158
$ws.Range("M83").Value = -3229.00014
$ws.Range("N83").Value = -30646.6158

$ws.Range("H88").Value = 2119.65
$ws.Range("I88").Value = 2164.6
$ws.Range("J88").Value = 2104.6667
$ws.Range("K88").Value = 2164.6
$ws.Range("L88").Value = 2104.6667
$ws.Range("M88").Value = -1758.6
$ws.Range("N88").Value = -2916.6667

$ws.Range("H91").Value = 2119.65
$ws.Range("I91").Value = 2164.6
$ws.Range("J91").Value = 2104.6667
$ws.Range("K91").Value = 2164.6
$ws.Range("L91").Value = 2104.6667
$ws.Range("M91").Value = -760.5999999999999
$ws.Range("N91").Value = -4912.6667

$ws.Range("H100").Value = 2510.2222
$ws.Range("I100").Value = 2450
$ws.Range("J100").Value = 2992
$ws.Range("K100").Value = 2450
$ws.Range("L100").Value = 2992
$ws.Range("M100").Value = -1909
$ws.Range("N100").Value = -4074

$ws.Range("H112").Value = 3103.9412
$ws.Range("J112").Value = 3210.5625
$ws.Range("L112").Value = 9631.6875
$ws.Range("N112").Value = -11847.6875

$ws.Range("H116").Value = 5000
$ws.Range("J116").Value = 5000
$ws.Range("L116").Value = 5000
$ws.Range("N116").Value = -11884

$ws.Range("H118").Value = 1715.3334
$ws.Range("I118").Value = 323.25
$ws.Range("K118").Value = 969.75
$ws.Range("M118").Value = 687.25

$ws.Range("H120").Value = 99999
$ws.Range("J120").Value = 99999
$ws.Range("L120").Value = 99999
$ws.Range("N120").Value = -109675

$ws.Range("H138").Value = 1577.3889
$ws.Range("I138").Value = 1405.4706
$ws.Range("J138").Value = 4500
$ws.Range("K138").Value = 4216.4118
$ws.Range("L138").Value = 13500
$ws.Range("M138").Value = 923.5882000000001
$ws.Range("N138").Value = -23780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2905.8
$ws.Range("I2").Value = 1560.6364
$ws.Range("J2").Value = 4549.8887
$ws.Range("K2").Value = 1560.6364
$ws.Range("L2").Value = 4549.8887
$ws.Range("M2").Value = -1447.6364
$ws.Range("N2").Value = -4775.8887

$ws.Range("H74").Value = 1196.2858
$ws.Range("I74").Value = 1265.6666
$ws.Range("J74").Value = 780
$ws.Range("K74").Value = 1265.6666
$ws.Range("L74").Value = 780
$ws.Range("M74").Value = -391.6666
$ws.Range("N74").Value = -2528

$ws.Range("H77").Value = 1196.2858
$ws.Range("I77").Value = 1265.6666
$ws.Range("J77").Value = 780
$ws.Range("K77").Value = 6328.333000000001
$ws.Range("L77").Value = 3900
$ws.Range("M77").Value = -1960.333000000001
$ws.Range("N77").Value = -12636

$ws.Range("H116").Value = 2905.8
$ws.Range("I116").Value = 1560.6364
$ws.Range("J116").Value = 4549.8887
$ws.Range("K116").Value = 1560.6364
$ws.Range("L116").Value = 4549.8887
$ws.Range("M116").Value = 733.3635999999999
$ws.Range("N116").Value = -9137.8887

$ws.Range("H132").Value = 402.46155
$ws.Range("I132").Value = 294.0909
$ws.Range("K132").Value = 882.2727
$ws.Range("M132").Value = 1647.7273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2905.8
$ws.Range("I3").Value = 1560.6364
$ws.Range("J3").Value = 4549.8887
$ws.Range("K3").Value = 1560.6364
$ws.Range("L3").Value = 4549.8887
$ws.Range("M3").Value = -1446.6364
$ws.Range("N3").Value = -4777.8887

$ws.Range("H20").Value = 3341.2856
$ws.Range("J20").Value = 4297.5713
$ws.Range("L20").Value = 4297.5713
$ws.Range("N20").Value = -4791.5713

$ws.Range("H82").Value = 17300
$ws.Range("I82").Value = 6400
$ws.Range("J82").Value = 50000
$ws.Range("K82").Value = 6400
$ws.Range("L82").Value = 50000
$ws.Range("M82").Value = -6017
$ws.Range("N82").Value = -50766

$ws.Range("H85").Value = 17300
$ws.Range("I85").Value = 6400
$ws.Range("J85").Value = 50000
$ws.Range("K85").Value = 6400
$ws.Range("L85").Value = 50000
$ws.Range("M85").Value = -5074
$ws.Range("N85").Value = -52652

$ws.Range("H86").Value = 9166.666999999999
$ws.Range("I86").Value = 1500
$ws.Range("J86").Value = 13000
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 13000
$ws.Range("M86").Value = -377
$ws.Range("N86").Value = -15246

$ws.Range("H89").Value = 9166.666999999999
$ws.Range("I89").Value = 1500
$ws.Range("J89").Value = 13000
$ws.Range("K89").Value = 7500
$ws.Range("L89").Value = 65000
$ws.Range("M89").Value = -1884
$ws.Range("N89").Value = -76232

$ws.Range("H105").Value = 3470.7222
$ws.Range("I105").Value = 3077.5715
$ws.Range("K105").Value = 3077.5715
$ws.Range("M105").Value = -1330.5715

$ws.Range("H134").Value = 2776.5
$ws.Range("I134").Value = 2776.5
$ws.Range("K134").Value = 8329.5
$ws.Range("M134").Value = -5794.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1385
$ws.Range("I58").Value = 965.0769
$ws.Range("K58").Value = 965.0769
$ws.Range("M58").Value = -762.0769

$ws.Range("H86").Value = 14341593
$ws.Range("I86").Value = 14341593
$ws.Range("K86").Value = 14341593
$ws.Range("M86").Value = -14340470

$ws.Range("H89").Value = 14341593
$ws.Range("I89").Value = 14341593
$ws.Range("K89").Value = 71707965
$ws.Range("M89").Value = -71702349

$ws.Range("H134").Value = 1465.8572
$ws.Range("I134").Value = 1117.7
$ws.Range("J134").Value = 2336.25
$ws.Range("K134").Value = 3353.1
$ws.Range("L134").Value = 7008.75
$ws.Range("M134").Value = -818.1000000000004
$ws.Range("N134").Value = -12078.75

$ws.Range("H136").Value = 1385
$ws.Range("I136").Value = 965.0769
$ws.Range("K136").Value = 2895.2307
$ws.Range("M136").Value = -345.2307000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 8000
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 8000
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 24000
$ws.Range("N104").Value = -29242
$ws.Range("M104").ClearContents()

$ws.Range("H132").Value = 2197
$ws.Range("J132").Value = 2496.25
$ws.Range("L132").Value = 22466.25
$ws.Range("N132").Value = -27526.25

$ws.Range("H139").Value = 4232.5
$ws.Range("I139").Value = 4232.5
$ws.Range("K139").Value = 12697.5
$ws.Range("M139").Value = -7557.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()

$ws.Range("H126").Value = 4149.75
$ws.Range("J126").Value = 4149.75
$ws.Range("L126").Value = 12449.25
$ws.Range("N126").Value = -17389.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1018.8
$ws.Range("J22").Value = 998
$ws.Range("L22").Value = 998
$ws.Range("N22").Value = -1588

$ws.Range("H27").Value = 1018.8
$ws.Range("J27").Value = 998
$ws.Range("L27").Value = 998
$ws.Range("N27").Value = -1212

$ws.Range("H46").Value = 1879.1578
$ws.Range("I46").Value = 1642.8572
$ws.Range("J46").Value = 2017
$ws.Range("K46").Value = 1642.8572
$ws.Range("L46").Value = 2017
$ws.Range("M46").Value = -1454.8572
$ws.Range("N46").Value = -2393

$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 2000
$ws.Range("K61").Value = 2000
$ws.Range("M61").Value = -1798

$ws.Range("H93").Value = 1713.1
$ws.Range("I93").Value = 1538.5
$ws.Range("J93").Value = 1829.5
$ws.Range("K93").Value = 1538.5
$ws.Range("L93").Value = 1829.5
$ws.Range("M93").Value = -290.5
$ws.Range("N93").Value = -4325.5

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 170

$ws.Range("H136").Value = 1912.8572
$ws.Range("I136").Value = 1906.1538
$ws.Range("K136").Value = 5718.4614
$ws.Range("M136").Value = -3168.4614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2679
$ws.Range("I81").Value = 3075
$ws.Range("J81").Value = 1095
$ws.Range("K81").Value = 6150
$ws.Range("L81").Value = 2190
$ws.Range("M81").Value = -5089
$ws.Range("N81").Value = -4312

$ws.Range("H84").Value = 2679
$ws.Range("I84").Value = 3075
$ws.Range("J84").Value = 1095
$ws.Range("K84").Value = 30750
$ws.Range("L84").Value = 10950
$ws.Range("M84").Value = -25446
$ws.Range("N84").Value = -21558

$ws.Range("H113").Value = 561.875
$ws.Range("I113").Value = 549
$ws.Range("K113").Value = 1647
$ws.Range("M113").Value = 523
